# perbaikan import dan excel mahasiswa
# Remove the "id_kategori" column (J) contents - header and data -
# which also drops the now-unused "id_kategori" shared string,
# and move the active selection to J6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the id_kategori header (J1) and its data value (J2).
$ws.Range("J1").ClearContents()
$ws.Range("J2").ClearContents()

# Update the active selection to match the saved view state.
$ws.Range("J6").Select()
